$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# habitat_typenumber values (column F) to add/update, per commit diff
$updates = @{
    5 = "6000"
    7 = "6233"
    10 = "9911"
    13 = "9915"
    17 = "6233"
    18 = "9990"
    27 = "9940"
    30 = "6000"
    33 = "9920"
    34 = "9920"
    35 = "6233"
    36 = "9940"
    41 = "6212"
    44 = "9990"
    53 = "9911"
    55 = "6000"
    60 = "9990"
    61 = "6402"
    69 = "9915"
    74 = "6233"
    75 = "7007"
    76 = "6402"
    77 = "9940"
    78 = "9911"
    79 = "9990"
    86 = "9920"
    87 = "9911"
    88 = "9990"
    96 = "6233"
    99 = "9940"
    101 = "9920"
    103 = "6233"
    104 = "6000"
    106 = "9990"
    108 = "9990"
    114 = "9911"
    116 = "6000"
    119 = "9920"
    120 = "6402"
    121 = "6233"
    125 = "6402"
    126 = "6402"
    130 = "6233"
    131 = "9940"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 6)   # column F = 6
    $cell.NumberFormat = "@"           # ensure text storage, matching t="inlineStr" source values
    $cell.Value = $updates[$row]
}
